$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 273.77274
$ws.Range("I28").Value = 172.46666
$ws.Range("J28").Value = 490.85715
$ws.Range("K28").Value = 172.46666
$ws.Range("L28").Value = 490.85715
$ws.Range("M28").Value = 312.53334
$ws.Range("N28").Value = -1460.85715
# Row 80
$ws.Range("H80").Value = 365.5357
$ws.Range("I80").Value = 346.66666
$ws.Range("J80").Value = 387.30768
$ws.Range("K80").Value = 1039.99998
$ws.Range("L80").Value = 1161.92304
$ws.Range("M80").Value = -41.99998000000005
$ws.Range("N80").Value = -3157.92304
# Row 83
$ws.Range("H83").Value = 365.5357
$ws.Range("I83").Value = 346.66666
$ws.Range("J83").Value = 387.30768
$ws.Range("K83").Value = 3119.99994
$ws.Range("L83").Value = 3485.76912
$ws.Range("M83").Value = 1872.00006
$ws.Range("N83").Value = -13469.76912
# Row 86
$ws.Range("H86").Value = 26500
$ws.Range("J86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("N86").Value = -4246
# Row 89
$ws.Range("H89").Value = 26500
$ws.Range("J89").Value = 2000
$ws.Range("L89").Value = 10000
$ws.Range("N89").Value = -21232
# Row 92
$ws.Range("H92").Value = 521.9666999999999
$ws.Range("I92").Value = 330.21054
$ws.Range("J92").Value = 853.1818
$ws.Range("K92").Value = 330.21054
$ws.Range("L92").Value = 853.1818
$ws.Range("M92").Value = 917.78946
$ws.Range("N92").Value = -3349.1818
# Row 106
$ws.Range("H106").Value = 38463036
$ws.Range("I106").Value = 62500908
$ws.Range("J106").Value = 2443
$ws.Range("K106").Value = 62500908
$ws.Range("L106").Value = 2443
$ws.Range("M106").Value = -62500277
$ws.Range("N106").Value = -3705
# Row 138
$ws.Range("H138").Value = 1750.1134
$ws.Range("I138").Value = 905.0238000000001
$ws.Range("J138").Value = 2395.4546
$ws.Range("K138").Value = 2715.0714
$ws.Range("L138").Value = 7186.3638
$ws.Range("M138").Value = 2424.9286
$ws.Range("N138").Value = -17466.3638

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 741.625
$ws.Range("I2").Value = 596.6
$ws.Range("J2").Value = 983.3333
$ws.Range("K2").Value = 596.6
$ws.Range("L2").Value = 983.3333
$ws.Range("M2").Value = -483.6
$ws.Range("N2").Value = -1209.3333
# Row 32
$ws.Range("H32").Value = 2805.74
$ws.Range("I32").Value = 2657.5
$ws.Range("J32").Value = 4139.9
$ws.Range("K32").Value = 2657.5
$ws.Range("L32").Value = 4139.9
$ws.Range("M32").Value = -2370.5
$ws.Range("N32").Value = -4713.9
# Row 97
$ws.Range("H97").Value = 936.1892
$ws.Range("I97").Value = 809.1923
$ws.Range("J97").Value = 1236.3636
$ws.Range("K97").Value = 809.1923
$ws.Range("L97").Value = 1236.3636
$ws.Range("M97").Value = -313.1923
$ws.Range("N97").Value = -2228.3636
# Row 107
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680
# Row 116
$ws.Range("H116").Value = 741.625
$ws.Range("I116").Value = 596.6
$ws.Range("J116").Value = 983.3333
$ws.Range("K116").Value = 596.6
$ws.Range("L116").Value = 983.3333
$ws.Range("M116").Value = 1697.4
$ws.Range("N116").Value = -5571.3333
# Row 132
$ws.Range("H132").Value = 31254742
$ws.Range("I132").Value = 40001270
$ws.Range("K132").Value = 120003810
$ws.Range("M132").Value = -120001280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 741.625
$ws.Range("I3").Value = 596.6
$ws.Range("J3").Value = 983.3333
$ws.Range("K3").Value = 596.6
$ws.Range("L3").Value = 983.3333
$ws.Range("M3").Value = -482.6
$ws.Range("N3").Value = -1211.3333
# Row 107
$ws.Range("H107").Value = 8772456
$ws.Range("I107").Value = 11364065
$ws.Range("J107").Value = 856.0769
$ws.Range("K107").Value = 11364065
$ws.Range("L107").Value = 856.0769
$ws.Range("M107").Value = -11362145
$ws.Range("N107").Value = -4696.0769
# Row 134
$ws.Range("H134").Value = 1918856.8
$ws.Range("I134").Value = 1059.1
$ws.Range("J134").Value = 6180629
$ws.Range("K134").Value = 3177.3
$ws.Range("L134").Value = 18541887
$ws.Range("M134").Value = -642.2999999999997
$ws.Range("N134").Value = -18546957

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1633.079
$ws.Range("I31").Value = 1132
$ws.Range("J31").Value = 2134.158
$ws.Range("K31").Value = 1132
$ws.Range("L31").Value = 2134.158
$ws.Range("M31").Value = -837
$ws.Range("N31").Value = -2724.158
# Row 34
$ws.Range("H34").Value = 1633.079
$ws.Range("I34").Value = 1132
$ws.Range("J34").Value = 2134.158
$ws.Range("K34").Value = 1132
$ws.Range("L34").Value = 2134.158
$ws.Range("M34").Value = -930
$ws.Range("N34").Value = -2538.158
# Row 122
$ws.Range("H122").Value = 16668484
$ws.Range("I122").Value = 22728032
$ws.Range("J122").Value = 4725
$ws.Range("K122").Value = 68184096
$ws.Range("L122").Value = 14175
$ws.Range("M122").Value = -68181646
$ws.Range("N122").Value = -19075

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 20
$ws.Range("H20").Value = 500900
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 5000000
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 15000000
$ws.Range("M20").Value = -2773
$ws.Range("N20").Value = -15000454
# Row 122
$ws.Range("H122").Value = 26321690
$ws.Range("I122").Value = 45454856
$ws.Range("J122").Value = 13587.375
$ws.Range("K122").Value = 409093704
$ws.Range("L122").Value = 122286.375
$ws.Range("M122").Value = -409091254
$ws.Range("N122").Value = -127186.375
# Row 131
$ws.Range("H131").Value = 770.15
$ws.Range("I131").Value = 464.66666
$ws.Range("J131").Value = 824.05884
$ws.Range("K131").Value = 1393.99998
$ws.Range("L131").Value = 2472.17652
$ws.Range("M131").Value = 3646.00002
$ws.Range("N131").Value = -12552.17652
# Row 137
$ws.Range("H137").Value = 31252234
$ws.Range("I137").Value = 35716336
$ws.Range("J137").Value = 3516.5
$ws.Range("K137").Value = 107149008
$ws.Range("L137").Value = 10549.5
$ws.Range("M137").Value = -107143908
$ws.Range("N137").Value = -20749.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1040.9474
$ws.Range("I102").Value = 998.7059
$ws.Range("K102").Value = 998.7059
$ws.Range("M102").Value = 623.2941

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 52637196
$ws.Range("I55").Value = 8046.5386
$ws.Range("J55").Value = 166667020
$ws.Range("K55").Value = 8046.5386
$ws.Range("L55").Value = 166667020
$ws.Range("M55").Value = -7873.5386
$ws.Range("N55").Value = -166667366
# Row 122
$ws.Range("H122").Value = 5538.4863
$ws.Range("I122").Value = 5800.6895
$ws.Range("J122").Value = 4588
$ws.Range("K122").Value = 17402.0685
$ws.Range("L122").Value = 13764
$ws.Range("M122").Value = -14952.0685
$ws.Range("N122").Value = -18664
